$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update temperature values (TemperatureC / TemperatureF) and the Summary text
# for each forecast row, per the new weather results.

$ws.Range("B2").Value = 47
$ws.Range("C2").Value = 116
$ws.Range("D2").Value = "Freezing"

$ws.Range("B3").Value = 29
$ws.Range("C3").Value = 84
$ws.Range("D3").Value = "Sweltering"

$ws.Range("B4").Value = 34
$ws.Range("C4").Value = 93
$ws.Range("D4").Value = "Scorching"

$ws.Range("B5").Value = -1
$ws.Range("C5").Value = 31
$ws.Range("D5").Value = "Balmy"

$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = "Mild"

# Slightly widen the Summary column to fit the new text
$ws.Columns.Item(4).ColumnWidth = 9.166666666666666
